# Estadisticos Matutinos 15 Oct
# - Recompute the "1er Parcial" (Estadisticos 1P) and "Estadisticos Final"
#   statistics for groups 1AM, 1BM and 3AEV (rows 2, 3 and 6) now that the
#   rescatable/failing students have been identified.
# - Fix the "Aprobados" (E) column on Estadisticos 2P for the same groups.
# - Populate the "Rescatables" sheet with the students who need to retake
#   the exam.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Estadisticos 1P — rows 2, 3, 6: Blancos/Reprobados/Aprobados/Por_Apro/Promedio
# ---------------------------------------------------------------------
$ws1P = $wb.Worksheets.Item("Estadisticos 1P")

$ws1P.Range("D2").Value = 0
$ws1P.Range("E2").Value = 10
$ws1P.Range("F2").Value = 24
$ws1P.Range("G2").Value = 70.59
$ws1P.Range("H2").Value = 7.6

$ws1P.Range("D3").Value = 0
$ws1P.Range("E3").Value = 10
$ws1P.Range("F3").Value = 21
$ws1P.Range("G3").Value = 67.74
$ws1P.Range("H3").Value = 7.3

$ws1P.Range("D6").Value = 0
$ws1P.Range("E6").Value = 12
$ws1P.Range("F6").Value = 21
$ws1P.Range("G6").Value = 63.64
$ws1P.Range("H6").Value = 6.6

# ---------------------------------------------------------------------
# Estadisticos 2P — rows 2, 3, 6: Aprobados (E) should equal Totales (D)
# ---------------------------------------------------------------------
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")

$ws2P.Range("E2").Value = 34
$ws2P.Range("E3").Value = 31
$ws2P.Range("E6").Value = 33

# ---------------------------------------------------------------------
# Estadisticos Final — rows 2, 3, 6: same recompute as Estadisticos 1P
# ---------------------------------------------------------------------
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")

$wsFinal.Range("D2").Value = 0
$wsFinal.Range("E2").Value = 10
$wsFinal.Range("F2").Value = 24
$wsFinal.Range("G2").Value = 70.59
$wsFinal.Range("H2").Value = 7.6

$wsFinal.Range("D3").Value = 0
$wsFinal.Range("E3").Value = 10
$wsFinal.Range("F3").Value = 21
$wsFinal.Range("G3").Value = 67.74
$wsFinal.Range("H3").Value = 7.3

$wsFinal.Range("D6").Value = 0
$wsFinal.Range("E6").Value = 12
$wsFinal.Range("F6").Value = 21
$wsFinal.Range("G6").Value = 63.64
$wsFinal.Range("H6").Value = 6.6

# ---------------------------------------------------------------------
# Rescatables — add the 4 students who need to retake the exam
# Columns: NC | Paterno | Materno | Nombres | Mat | Grupo | Reprobadas
# ---------------------------------------------------------------------
$wsResc = $wb.Worksheets.Item("Rescatables")

$wsResc.Cells.Item(2, 1).Value = 21330051920162
$wsResc.Cells.Item(2, 2).Value = "IXMATLAHUA"
$wsResc.Cells.Item(2, 3).Value = "HERNANDEZ"
$wsResc.Cells.Item(2, 4).Value = "CRISTIAN FERNANDO"
$wsResc.Cells.Item(2, 5).Value = "ÁLGEBRA"
$wsResc.Cells.Item(2, 6).Value = "1AM"
$wsResc.Cells.Item(2, 7).Value = 6

$wsResc.Cells.Item(3, 1).Value = 20330051920284
$wsResc.Cells.Item(3, 2).Value = "ANTONIO"
$wsResc.Cells.Item(3, 3).Value = "GARCIA"
$wsResc.Cells.Item(3, 4).Value = "ISRAEL"
$wsResc.Cells.Item(3, 5).Value = "GEOMETRÍA ANALÍTICA"
$wsResc.Cells.Item(3, 6).Value = "3BLCM"
$wsResc.Cells.Item(3, 7).Value = 6

$wsResc.Cells.Item(4, 1).Value = 21330051920207
$wsResc.Cells.Item(4, 2).Value = "ROBLES"
$wsResc.Cells.Item(4, 3).Value = "SANCHEZ"
$wsResc.Cells.Item(4, 4).Value = "SINAI ANTONIO"
$wsResc.Cells.Item(4, 5).Value = "ÁLGEBRA"
$wsResc.Cells.Item(4, 6).Value = "1BM"
$wsResc.Cells.Item(4, 7).Value = 6

$wsResc.Cells.Item(5, 1).Value = 21330051920213
$wsResc.Cells.Item(5, 2).Value = "VENTURA"
$wsResc.Cells.Item(5, 3).Value = "ROSALES"
$wsResc.Cells.Item(5, 4).Value = "GUSTAVO"
$wsResc.Cells.Item(5, 5).Value = "ÁLGEBRA"
$wsResc.Cells.Item(5, 6).Value = "1BM"
$wsResc.Cells.Item(5, 7).Value = 6

Write-Output "edit applied"
